# "error solve ifrs list" - replace the (incorrectly scaled) financial figures
# in the "company_list" sheet with the corrected values, row by row.
# Some rows also drop the AD/AH columns entirely (no longer reported),
# which is done with ClearContents() so the cells disappear from the XML.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 824
$ws.Range("E2").Value = 181
$ws.Range("F2").Value = 169
$ws.Range("G2").Value = 194
$ws.Range("H2").Value = 171
$ws.Range("I2").Value = 170
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1138
$ws.Range("L2").Value = 728
$ws.Range("M2").Value = 410
$ws.Range("N2").Value = 408
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 62
$ws.Range("Q2").Value = 101
$ws.Range("R2").Value = -214
$ws.Range("S2").Value = 130
$ws.Range("T2").Value = 222
$ws.Range("U2").Value = -121
$ws.Range("V2").Value = 542
$ws.Range("W2").Value = 21.98
$ws.Range("X2").Value = 20.74
$ws.Range("Y2").Value = 41.17
$ws.Range("Z2").Value = 17.14
$ws.Range("AA2").Value = 177.35
$ws.Range("AB2").Value = 1190.3
$ws.Range("AC2").Value = 1542
$ws.Range("AE2").Value = 3687
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 70
$ws.Range("AI2").Value = 2.5
$ws.Range("AJ2").Value = 11068830
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 955
$ws.Range("E3").Value = 168
$ws.Range("F3").Value = 168
$ws.Range("G3").Value = 174
$ws.Range("H3").Value = 145
$ws.Range("I3").Value = 145
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1506
$ws.Range("L3").Value = 1019
$ws.Range("M3").Value = 488
$ws.Range("N3").Value = 485
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 62
$ws.Range("Q3").Value = 150
$ws.Range("R3").Value = -283
$ws.Range("S3").Value = 149
$ws.Range("T3").Value = 279
$ws.Range("U3").Value = -129
$ws.Range("V3").Value = 698
$ws.Range("W3").Value = 17.57
$ws.Range("X3").Value = 15.18
$ws.Range("Y3").Value = 32.42
$ws.Range("Z3").Value = 10.96
$ws.Range("AA3").Value = 208.92
$ws.Range("AB3").Value = 1317.37
$ws.Range("AC3").Value = 1308
$ws.Range("AE3").Value = 4384
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 100
$ws.Range("AI3").Value = 4.2
$ws.Range("AJ3").Value = 11068830
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 1200
$ws.Range("E4").Value = 285
$ws.Range("F4").Value = 285
$ws.Range("G4").Value = 264
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 198
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 1874
$ws.Range("L4").Value = 1179
$ws.Range("M4").Value = 694
$ws.Range("N4").Value = 690
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 62
$ws.Range("Q4").Value = 176
$ws.Range("R4").Value = -191
$ws.Range("S4").Value = 52
$ws.Range("T4").Value = 182
$ws.Range("U4").Value = -6
$ws.Range("V4").Value = 779
$ws.Range("W4").Value = 23.77
$ws.Range("X4").Value = 16.66
$ws.Range("Y4").Value = 33.68
$ws.Range("Z4").Value = 11.83
$ws.Range("AA4").Value = 169.81
$ws.Range("AB4").Value = 1658.2
$ws.Range("AC4").Value = 1788
$ws.Range("AE4").Value = 6234
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 100
$ws.Range("AI4").Value = 3.07
$ws.Range("AJ4").Value = 11068830
$ws.Range("AD4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 1506
$ws.Range("E5").Value = 411
$ws.Range("F5").Value = 411
$ws.Range("G5").Value = 359
$ws.Range("H5").Value = 303
$ws.Range("I5").Value = 301
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 3109
$ws.Range("L5").Value = 1484
$ws.Range("M5").Value = 1625
$ws.Range("N5").Value = 1619
$ws.Range("O5").Value = 6
$ws.Range("P5").Value = 62
$ws.Range("Q5").Value = 217
$ws.Range("R5").Value = -739
$ws.Range("S5").Value = 964
$ws.Range("T5").Value = 723
$ws.Range("U5").Value = -507
$ws.Range("V5").Value = 974
$ws.Range("W5").Value = 27.27
$ws.Range("X5").Value = 20.11
$ws.Range("Y5").Value = 26.05
$ws.Range("Z5").Value = 12.16
$ws.Range("AA5").Value = 91.34999999999999
$ws.Range("AB5").Value = 2854.73
$ws.Range("AC5").Value = 2716
$ws.Range("AD5").Value = 22.94
$ws.Range("AE5").Value = 18768
$ws.Range("AF5").Value = 3.32
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 0.24
$ws.Range("AI5").Value = 4.3
$ws.Range("AJ5").Value = 11068830

# Row 6
$ws.Range("D6").Value = 1863
$ws.Range("E6").Value = 420
$ws.Range("F6").Value = 420
$ws.Range("G6").Value = 420
$ws.Range("H6").Value = 343
$ws.Range("I6").Value = 343
$ws.Range("K6").Value = 4368
$ws.Range("L6").Value = 2421
$ws.Range("M6").Value = 1947
$ws.Range("N6").Value = 1941
$ws.Range("P6").Value = 62
$ws.Range("Q6").Value = -6
$ws.Range("R6").Value = -1118
$ws.Range("S6").Value = 802
$ws.Range("T6").Value = 1074
$ws.Range("U6").Value = -1080
$ws.Range("V6").Value = 1933
$ws.Range("W6").Value = 22.52
$ws.Range("X6").Value = 18.39
$ws.Range("Y6").Value = 19.25
$ws.Range("Z6").Value = 9.16
$ws.Range("AA6").Value = 124.32
$ws.Range("AB6").Value = 3378.88
$ws.Range("AC6").Value = 3095
$ws.Range("AD6").Value = 18.58
$ws.Range("AE6").Value = 22506
$ws.Range("AF6").Value = 2.55
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 0.26
$ws.Range("AI6").Value = 3.78
$ws.Range("AJ6").Value = 11068830

# Row 7
$ws.Range("D7").Value = 2551
$ws.Range("E7").Value = 504
$ws.Range("G7").Value = 438
$ws.Range("H7").Value = 325
$ws.Range("I7").Value = 327
$ws.Range("K7").Value = 5144
$ws.Range("L7").Value = 2909
$ws.Range("M7").Value = 2235
$ws.Range("N7").Value = 2233
$ws.Range("P7").Value = 60
$ws.Range("Q7").Value = 231
$ws.Range("R7").Value = -539
$ws.Range("S7").Value = 380
$ws.Range("T7").Value = 497
$ws.Range("U7").Value = -416
$ws.Range("W7").Value = 19.74
$ws.Range("X7").Value = 12.74
$ws.Range("Y7").Value = 15.67
$ws.Range("Z7").Value = 6.83
$ws.Range("AA7").Value = 130.16
$ws.Range("AC7").Value = 2954
$ws.Range("AD7").Value = 16.47
$ws.Range("AE7").Value = 25897
$ws.Range("AF7").Value = 1.88
$ws.Range("AG7").Value = 150
$ws.Range("AH7").Value = 0.31
$ws.Range("AI7").Value = 5.09

# Row 8
$ws.Range("D8").Value = 3092
$ws.Range("E8").Value = 624
$ws.Range("G8").Value = 596
$ws.Range("H8").Value = 465
$ws.Range("I8").Value = 469
$ws.Range("K8").Value = 5672
$ws.Range("L8").Value = 2998
$ws.Range("M8").Value = 2674
$ws.Range("N8").Value = 2683
$ws.Range("P8").Value = 60
$ws.Range("Q8").Value = 440
$ws.Range("R8").Value = -312
$ws.Range("S8").Value = -4
$ws.Range("T8").Value = 273
$ws.Range("U8").Value = 63
$ws.Range("W8").Value = 20.18
$ws.Range("X8").Value = 15.03
$ws.Range("Y8").Value = 19.08
$ws.Range("Z8").Value = 8.59
$ws.Range("AA8").Value = 112.14
$ws.Range("AC8").Value = 4237
$ws.Range("AD8").Value = 11.48
$ws.Range("AE8").Value = 31115
$ws.Range("AF8").Value = 1.56
$ws.Range("AG8").Value = 112
$ws.Range("AH8").Value = 0.23
$ws.Range("AI8").Value = 2.66

# Row 9
$ws.Range("D9").Value = 3731
$ws.Range("E9").Value = 790
$ws.Range("G9").Value = 765
$ws.Range("H9").Value = 599
$ws.Range("I9").Value = 601
$ws.Range("K9").Value = 6400
$ws.Range("L9").Value = 3141
$ws.Range("M9").Value = 3258
$ws.Range("N9").Value = 3287
$ws.Range("P9").Value = 60
$ws.Range("Q9").Value = 493
$ws.Range("R9").Value = -375
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 265
$ws.Range("U9").Value = 127
$ws.Range("W9").Value = 21.19
$ws.Range("X9").Value = 16.05
$ws.Range("Y9").Value = 20.14
$ws.Range("Z9").Value = 9.92
$ws.Range("AA9").Value = 96.39
$ws.Range("AC9").Value = 5432
$ws.Range("AD9").Value = 8.960000000000001
$ws.Range("AE9").Value = 38111
$ws.Range("AF9").Value = 1.28
$ws.Range("AG9").Value = 112
$ws.Range("AH9").Value = 0.23
$ws.Range("AI9").Value = 2.07
